# 1.3.6 inserção linha vsho11
#
# - Adds a new column P ("Preço Médio Pago") with a header matching the
#   existing header style, plus literal "avg price paid" figures for the
#   FII rows (7-14).
# - Inserts a new row at position 15 for the ticker "VSHO11" (all zeros),
#   pushing "Total Div Ações" and everything below it down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column P: header + style copied from the existing header (O1) ---
$ws.Range("P1").Value = "Preço Médio Pago"
$ws.Range("O1").Copy()
$ws.Range("P1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Literal "Preço Médio Pago" values for the FII rows (7-14) ---
$ws.Cells.Item(7, 16).Value = 205.2
$ws.Cells.Item(8, 16).Value = 3351.6
$ws.Cells.Item(9, 16).Value = 3227.14
$ws.Cells.Item(10, 16).Value = 3219.44
$ws.Cells.Item(11, 16).Value = 2995.3
$ws.Cells.Item(12, 16).Value = 3043.3
$ws.Cells.Item(13, 16).Value = 3042.39
$ws.Cells.Item(14, 16).Value = 3018.54

# --- Insert the new VSHO11 row at 15, shifting everything below down ---
$ws.Rows(15).Insert()

$ws.Cells.Item(15, 1).Value = "VSHO11"
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 0
$ws.Cells.Item(15, 16).Value = 0
